# transition rule 5 and 10 mi rad updates to script and all outputs
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Means"
# ---------------------------------------------------------------------------
$means = $wb.Worksheets.Item("Means")

# New headers for the 5-mile / 10-mile radius columns
$means.Range("F1").Value = "Within 5 miles of HFC production facility"
$means.Range("G1").Value = "Within 10 miles of HFC production facility"

# New column F (Within 5 miles) values, rows 2-10
$means.Range("F2").Value = 60
$means.Range("F3").Value = 30
$means.Range("F4").Value = 9.7
$means.Range("F5").Value = 11
$means.Range("F6").Value = 48
$means.Range("F7").Value = 12
$means.Range("F8").Value = 6.6
$means.Range("F9").Value = 40
$means.Range("F10").Value = 0.5

# New column G (Within 10 miles) values, rows 2-10
$means.Range("G2").Value = 72
$means.Range("G3").Value = 19
$means.Range("G4").Value = 9
$means.Range("G5").Value = 9.7
$means.Range("G6").Value = 57
$means.Range("G7").Value = 10
$means.Range("G8").Value = 5.7
$means.Range("G9").Value = 36
$means.Range("G10").Value = 0.45

# Updated Total Cancer Risk (row 9) and Total Respiratory (row 10) values
# for the existing National Average / State Average / 1 mile / 3 mile columns
$means.Range("B9").Value = 29
$means.Range("C9").Value = 34
$means.Range("D9").Value = 33
$means.Range("E9").Value = 42

$means.Range("B10").Value = 0.37
$means.Range("C10").Value = 0.47
$means.Range("D10").Value = 0.43
$means.Range("E10").Value = 0.52

# ---------------------------------------------------------------------------
# Sheet "Standard Deviations"
# ---------------------------------------------------------------------------
$sd = $wb.Worksheets.Item("Standard Deviations")

# New headers for the 5-mile / 10-mile radius SD columns
$sd.Range("F1").Value = "Within 5 mile of HFC production facility SD"
$sd.Range("G1").Value = "Within 10 mile of HFC production facility SD"

# New column F (Within 5 miles SD) values, rows 2-10
$sd.Range("F2").Value = 23
$sd.Range("F3").Value = 26
$sd.Range("F4").Value = 7.9
$sd.Range("F5").Value = 12
$sd.Range("F6").Value = 22
$sd.Range("F7").Value = 10
$sd.Range("F8").Value = 11
$sd.Range("F9").Value = 8.3
$sd.Range("F10").Value = 0.083

# New column G (Within 10 miles SD) values, rows 2-10
$sd.Range("G2").Value = 21
$sd.Range("G3").Value = 21
$sd.Range("G4").Value = 10
$sd.Range("G5").Value = 13
$sd.Range("G6").Value = 23
$sd.Range("G7").Value = 11
$sd.Range("G8").Value = 7.3
$sd.Range("G9").Value = 6.5
$sd.Range("G10").Value = 0.066

# Updated Total Cancer Risk (row 9) and Total Respiratory (row 10) values
# for the existing National Average / State Average / 1 mile / 3 mile columns
$sd.Range("B9").Value = 10
$sd.Range("C9").Value = 5.3
$sd.Range("D9").Value = 7.1
$sd.Range("E9").Value = 11

$sd.Range("B10").Value = 0.14
$sd.Range("C10").Value = 0.07
$sd.Range("D10").Value = 0.071
$sd.Range("E10").Value = 0.11

Write-Output "done"
